$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 51970749820
$ws.Range("A2").Value = 51946361875

$ws.Range("A4").Select()
